# [TEST SCRAPE] updated files from azure vm
#
# 1. Remove the stray empty INNING_NUMBER (column B) cells on the
#    "ODI Batting" sheet for rows that have no value (2, 6, 9, 12, 13, 16).
# 2. Add a new "ODI Batting Extra" worksheet at the end of the workbook
#    with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#    PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clear the empty inline-string cells in column B of ODI Batting
# ---------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$emptyBRows = @(2, 6, 9, 12, 13, 16)
foreach ($r in $emptyBRows) {
    $battingSheet.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------
# Step 2: add the new "ODI Batting Extra" sheet after "ODI Bowling"
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$extraSheet = $wb.Worksheets.Add([System.Type]::Missing, $bowlingSheet)
$extraSheet.Name = "ODI Batting Extra"

# Match the page margins used on the rest of the sheets in this workbook
# (0.75in left/right, 1in top/bottom, 0.5in header/footer).
$extraSheet.PageSetup.LeftMargin = 54
$extraSheet.PageSetup.RightMargin = 54
$extraSheet.PageSetup.TopMargin = 72
$extraSheet.PageSetup.BottomMargin = 72
$extraSheet.PageSetup.HeaderMargin = 36
$extraSheet.PageSetup.FooterMargin = 36

# Helper: force a cell to store its value as literal text, matching the
# inline-string cells used throughout this workbook (numeric-looking
# strings like "4115" or "5.86%" would otherwise be auto-detected as
# numbers/percentages by Excel).
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $extraSheet.Cells.Item(1, 1) "MATCH_CODE"
Set-TextValue $extraSheet.Cells.Item(1, 2) "BATTING_POSITION"
Set-TextValue $extraSheet.Cells.Item(1, 3) "NUM_4"
Set-TextValue $extraSheet.Cells.Item(1, 4) "NUM_6"
Set-TextValue $extraSheet.Cells.Item(1, 5) "PERCENT_RUNS_OF_TOTAL"
Set-TextValue $extraSheet.Cells.Item(1, 6) "MAN_OF_MATCH"

# Copy the header formatting (bold, centered, bordered) from the
# existing "ODI Batting" header row so the new header row matches the
# rest of the workbook's look. This must happen *after* the text is
# written, since Copy() overwrites both value and format, and writing
# text afterwards (Set-TextValue) would reset the style back to Normal.
$headerSrc = $battingSheet.Range("A1").Resize(1, 6)
$headerDst = $extraSheet.Range("A1").Resize(1, 6)
$headerSrc.Copy($headerDst)
$extraSheet.Cells.Item(1, 1).Value = "MATCH_CODE"
$extraSheet.Cells.Item(1, 2).Value = "BATTING_POSITION"
$extraSheet.Cells.Item(1, 3).Value = "NUM_4"
$extraSheet.Cells.Item(1, 4).Value = "NUM_6"
$extraSheet.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$extraSheet.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6,
# PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is numeric when present, everything else is text.
$rows = @(
    @{ Row = 2;  Code = "4115"; Pos = 11;  N4 = $null; N6 = $null; Pct = $null;     Mom = "NO"  },
    @{ Row = 3;  Code = "4167"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  },
    @{ Row = 4;  Code = "4168"; Pos = 9;   N4 = "0";   N6 = "0";   Pct = "5.86%";   Mom = "NO"  },
    @{ Row = 5;  Code = "4169"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  },
    @{ Row = 6;  Code = "4234"; Pos = 8;   N4 = $null; N6 = $null; Pct = $null;     Mom = "YES" },
    @{ Row = 7;  Code = "4235"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  },
    @{ Row = 8;  Code = "4236"; Pos = 8;   N4 = "1";   N6 = "0";   Pct = "6.96%";   Mom = "NO"  },
    @{ Row = 9;  Code = "4266"; Pos = 8;   N4 = $null; N6 = $null; Pct = $null;     Mom = "NO"  },
    @{ Row = 10; Code = "4268"; Pos = 8;   N4 = "0";   N6 = "0";   Pct = $null;     Mom = "NO"  },
    @{ Row = 11; Code = "4270"; Pos = 8;   N4 = "3";   N6 = "0";   Pct = "10.66%";  Mom = "NO"  },
    @{ Row = 12; Code = "4273"; Pos = 10;  N4 = $null; N6 = $null; Pct = $null;     Mom = "NO"  },
    @{ Row = 13; Code = "4274"; Pos = 10;  N4 = $null; N6 = $null; Pct = $null;     Mom = "NO"  },
    @{ Row = 14; Code = "4421"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  },
    @{ Row = 15; Code = "4594"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  },
    @{ Row = 16; Code = "4600"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;   Mom = "NO"  }
)

foreach ($item in $rows) {
    $r = $item.Row

    Set-TextValue $extraSheet.Cells.Item($r, 1) $item.Code

    if ($null -ne $item.Pos) {
        $extraSheet.Cells.Item($r, 2).Value = $item.Pos
    }

    if ($null -ne $item.N4) {
        Set-TextValue $extraSheet.Cells.Item($r, 3) $item.N4
    }

    if ($null -ne $item.N6) {
        Set-TextValue $extraSheet.Cells.Item($r, 4) $item.N6
    }

    if ($null -ne $item.Pct) {
        Set-TextValue $extraSheet.Cells.Item($r, 5) $item.Pct
    }

    Set-TextValue $extraSheet.Cells.Item($r, 6) $item.Mom
}
